$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 13
$years = 2010..2021
$K = @(5106,4543,5138,5869,5410,7066,4644,4919,8443,7991,5888,10665)
for ($i = 0; $i -lt 12; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 10).Value = $years[$i]
    $ws.Cells.Item($row, 11).Value = $K[$i]
}

$chartObj = $ws.Shapes.AddChart2(201, 4)
$chart = $chartObj.Chart

$ser1 = $chart.SeriesCollection.Item(1)
try {
    $ser1.Name = "=Sumary!`$K`$2"
    Write-Output "Name set ok"
} catch {
    Write-Output "Name set FAILED: $_"
}

try {
    $ser1.XValues = $ws.Range("J3:J14")
    Write-Output "XValues set ok"
} catch {
    Write-Output "XValues set FAILED: $_"
}

try {
    $ser1.Values = $ws.Range("K3:K14")
    Write-Output "Values set ok"
} catch {
    Write-Output "Values set FAILED: $_"
}

Write-Output "done"
